# Add print options UI
# Append new usage-log rows (36-42) to the "Наличные" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Наличные")

$rows = @(
    @(7676096317, "M-Банкинг чек-4294968802.pdf", 1, 0.2, "2025-06-28 14:22:42"),
    @(7676096317, "ПЕЧАТЬ11.docx", 2, 0.4, "2025-06-28 14:23:06"),
    @(7676096317, "M-Банкинг чек-4294968802.pdf", 1, 0.2, "2025-06-28 14:23:22"),
    @(7676096317, "M-Банкинг чек-4294968802.pdf", 1, 0.2, "2025-06-28 14:23:34"),
    @(7676096317, "M-Банкинг чек-4294968802.pdf", 1, 0.2, "2025-06-28 14:25:23"),
    @(7676096317, "357d1a45e0e15379ea555a7e8964ca76.pdf", 1, 0.2, "2025-06-28 14:29:31"),
    @(7676096317, "M-Банкинг чек.docx", 1, 0.2, "2025-06-28 14:33:43")
)

$startRow = 36
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
